$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.8
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 2.5
$ws.Range("Q2").Value = 1.74
$ws.Range("R2").Value = 2.04
$ws.Range("Z2").Value = 1.63
$ws.Range("AD2").Value = 13
$ws.Range("AI2").Value = 21
$ws.Range("AJ2").Value = 81
$ws.Range("AM2").Value = 23
$ws.Range("G3").Value = 1.22
$ws.Range("H3").Value = 7
$ws.Range("J3").Value = 1.47
$ws.Range("O3").Value = 1.1
$ws.Range("P3").Value = 7
$ws.Range("S3").Value = 1.33
$ws.Range("T3").Value = 3.4
$ws.Range("Y3").Value = 1.8
$ws.Range("Z3").Value = 1.95
$ws.Range("AB3").Value = 8
$ws.Range("AC3").Value = 10
$ws.Range("AF3").Value = 23
$ws.Range("AL3").Value = 34
$ws.Range("AN3").Value = 29
$ws.Range("AO3").Value = 126
$ws.Range("G4").Value = 1.8
$ws.Range("I4").Value = 5.25
$ws.Range("J4").Value = 2.6
$ws.Range("K4").Value = 1.87
$ws.Range("L4").Value = 6
$ws.Range("O4").Value = 1.57
$ws.Range("P4").Value = 2.25
$ws.Range("Q4").Value = 2.1
$ws.Range("R4").Value = 1.78
$ws.Range("U4").Value = 6
$ws.Range("V4").Value = 1.13
$ws.Range("W4").Value = 1.62
$ws.Range("X4").Value = 2.2
$ws.Range("Y4").Value = 2.5
$ws.Range("Z4").Value = 1.5
$ws.Range("AA4").Value = 4.75
$ws.Range("AD4").Value = 13
$ws.Range("AI4").Value = 23
$ws.Range("AJ4").Value = 101
$ws.Range("AM4").Value = 23
$ws.Range("AN4").Value = 19
$ws.Range("AP4").Value = 51
$ws.Range("AQ4").Value = 67
$ws.Range("G5").Value = 2.8
$ws.Range("H5").Value = 2.7
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 3.75
$ws.Range("K5").Value = 1.77
$ws.Range("M5").Value = 1.17
$ws.Range("N5").Value = 5
$ws.Range("O5").Value = 1.62
$ws.Range("P5").Value = 2.2
$ws.Range("S5").Value = 3.1
$ws.Range("T5").Value = 1.36
$ws.Range("U5").Value = 6.5
$ws.Range("V5").Value = 1.11
$ws.Range("AB5").Value = 12
$ws.Range("AC5").Value = 12
$ws.Range("AD5").Value = 29
$ws.Range("AF5").Value = 51
$ws.Range("AG5").Value = 5
$ws.Range("AI5").Value = 21
$ws.Range("AL5").Value = 6.5
$ws.Range("AR5").Value = 5.2
$ws.Range("AS5").Value = 1.16
$ws.Range("M6").Value = 1.17
$ws.Range("N6").Value = 5
$ws.Range("U6").Value = 8
$ws.Range("V6").Value = 1.08
$ws.Range("W6").Value = 1.8
$ws.Range("X6").Value = 2
$ws.Range("AS6").Value = 1.11
$ws.Range("G7").Value = 2.3
$ws.Range("I7").Value = 3.7
$ws.Range("J7").Value = 3.25
$ws.Range("N7").Value = 5
$ws.Range("Q7").Value = 2.6
$ws.Range("R7").Value = 1.5
$ws.Range("AA7").Value = 5
$ws.Range("AB7").Value = 9
$ws.Range("AD7").Value = 23
$ws.Range("AE7").Value = 29
$ws.Range("AL7").Value = 6.5
$ws.Range("AM7").Value = 15
$ws.Range("AQ7").Value = 51
$ws.Range("AR7").Value = 6.6
$ws.Range("G8").Value = 1.55
$ws.Range("I8").Value = 6.5
$ws.Range("J8").Value = 2.2
$ws.Range("M8").Value = 1.1
$ws.Range("N8").Value = 7
$ws.Range("Q8").Value = 1.83
$ws.Range("R8").Value = 2.03
$ws.Range("S8").Value = 2.4
$ws.Range("T8").Value = 1.53
$ws.Range("Y8").Value = 2.5
$ws.Range("Z8").Value = 1.5
$ws.Range("AG8").Value = 7
$ws.Range("AI8").Value = 26
$ws.Range("AL8").Value = 12
$ws.Range("AO8").Value = 81
$ws.Range("K9").Value = 1.8
$ws.Range("G15").Value = 2.2
$ws.Range("I15").Value = 3.2
$ws.Range("W15").Value = 1.37
$ws.Range("AA15").Value = 8
$ws.Range("AB15").Value = 11
$ws.Range("AC15").Value = 9.5
$ws.Range("AE15").Value = 19
$ws.Range("AG15").Value = 10
$ws.Range("AM15").Value = 15
$ws.Range("O16").Value = 1.3
$ws.Range("P16").Value = 3.4
$ws.Range("S16").Value = 2
$ws.Range("T16").Value = 1.85
$ws.Range("W16").Value = 1.37
$ws.Range("W17").Value = 1.58
$ws.Range("W18").Value = 1.58
$ws.Range("W20").Value = 1.3
$ws.Range("S24").Value = 1.87
$ws.Range("T24").Value = 1.87
$ws.Range("AP24").Value = 51
$ws.Range("N25").Value = 9
$ws.Range("S25").Value = 2.2
$ws.Range("T25").Value = 1.65
$ws.Range("U25").Value = 4
$ws.Range("V25").Value = 1.22
$ws.Range("Y25").Value = 1.8
$ws.Range("Z25").Value = 1.8
$ws.Range("U28").Value = 2.38
$ws.Range("V28").Value = 1.53
$ws.Range("AR28").Value = 1.9
$ws.Range("AS28").Value = 1.95
$ws.Range("H29").Value = 5.75
$ws.Range("I29").Value = 5.75
$ws.Range("J29").Value = 1.8
$ws.Range("Y29").Value = 1.36
$ws.Range("Z29").Value = 3
$ws.Range("AC29").Value = 10
$ws.Range("AK29").Value = 67
$ws.Range("AL29").Value = 34
$ws.Range("AP29").Value = 41
$ws.Range("G30").Value = 3.2
$ws.Range("I30").Value = 2.1
$ws.Range("J30").Value = 3.75
$ws.Range("W30").Value = 1.36
$ws.Range("X30").Value = 3
$ws.Range("AC30").Value = 12
$ws.Range("AO30").Value = 19
$ws.Range("G31").Value = 2.55
$ws.Range("H31").Value = 2.9
$ws.Range("I31").Value = 3
$ws.Range("J31").Value = 3.2
$ws.Range("L31").Value = 3.6
$ws.Range("S31").Value = 2.1
$ws.Range("T31").Value = 1.7
$ws.Range("Y31").Value = 1.8
$ws.Range("Z31").Value = 1.91
$ws.Range("AA31").Value = 8
$ws.Range("AB31").Value = 12
$ws.Range("AC31").Value = 10
$ws.Range("AD31").Value = 23
$ws.Range("AE31").Value = 21
$ws.Range("AL31").Value = 9
$ws.Range("AM31").Value = 15
$ws.Range("AS31").Value = 1.37
$ws.Range("G34").Value = 2.25
$ws.Range("I34").Value = 3.4
$ws.Range("J34").Value = 3.1
$ws.Range("Q34").Value = 1.77
$ws.Range("R34").Value = 1.97
$ws.Range("Y34").Value = 2.05
$ws.Range("Z34").Value = 1.7
$ws.Range("AA34").Value = 6
$ws.Range("AB34").Value = 9.5
$ws.Range("AD34").Value = 21
$ws.Range("AK34").Value = 501
$ws.Range("AN34").Value = 13
$ws.Range("AO34").Value = 41
$ws.Range("AP34").Value = 34
$ws.Range("G39").Value = 1.91
$ws.Range("H39").Value = 3.4
$ws.Range("I39").Value = 4
$ws.Range("J39").Value = 2.6
$ws.Range("K39").Value = 2.2
$ws.Range("L39").Value = 4.33
$ws.Range("O39").Value = 1.29
$ws.Range("P39").Value = 3.5
$ws.Range("S39").Value = 1.95
$ws.Range("T39").Value = 1.9
$ws.Range("Y39").Value = 1.69
$ws.Range("AA39").Value = 8
$ws.Range("AE39").Value = 15
$ws.Range("AG39").Value = 10
$ws.Range("AL39").Value = 12
$ws.Range("AM39").Value = 21
